$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 5.616015666666666
$ws.Range("N2").Value = 16.848047
$ws.Range("O2").Value = 0.2860808099623356
$ws.Range("P2").Value = 0.2860808099623357
$ws.Range("Q2").Value = 0.8617888360813333
$ws.Range("R2").Value = 7.756099524731999
$ws.Range("S2").Value = 0.2860808099623356
$ws.Range("T2").Value = 0.2860808099623357

# Row 3
$ws.Range("M3").Value = 8.435525999999999
$ws.Range("O3").Value = 0.4297071542841152
$ws.Range("P3").Value = 0.4297071542841153
$ws.Range("S3").Value = 0.4297071542841152
$ws.Range("T3").Value = 0.4297071542841153

# Row 4
$ws.Range("M4").Value = 2.036951
$ws.Range("N4").Value = 6.110853000000001
$ws.Range("O4").Value = 0.1037626364528048
$ws.Range("P4").Value = 0.1037626364528048
$ws.Range("Q4").Value = 0.312574204852
$ws.Range("R4").Value = 2.813167843668
$ws.Range("S4").Value = 0.1037626364528048
$ws.Range("T4").Value = 0.1037626364528048

# Row 5
$ws.Range("M5").Value = 3.542379
$ws.Range("N5").Value = 10.627137
$ws.Range("O5").Value = 0.1804493993007442
$ws.Range("P5").Value = 0.1804493993007443
$ws.Range("Q5").Value = 0.543585142308
$ws.Range("R5").Value = 4.892266280772001
$ws.Range("S5").Value = 0.1804493993007442
$ws.Range("T5").Value = 0.1804493993007443
